# Add the new "Before" worksheet, placed immediately after the existing
# "CellListener" sheet (so it becomes sheet index 2 / sheetId 2).
$wb = $excel.ActiveWorkbook
$firstSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $firstSheet)
$newSheet.Name = "Before"

# Seed content: a label, a placeholder the CellListener sample will
# overwrite at runtime, and an explanatory note.
$newSheet.Range("A1").Value = "State Name:"
$newSheet.Range("B1").Value = "Anything here; CellListener will replace!"
$newSheet.Range("B2").Value = 'The CellListener will replace the above content with ${california.name}'

# Match the template's column widths (18.42578125 / 36.7109375 "characters").
# This host's ColumnWidth setter only lands on 1/6-character increments, so
# these inputs are chosen to round to the closest representable width.
$newSheet.Columns.Item(1).ColumnWidth = 17.666666666666668
$newSheet.Columns.Item(2).ColumnWidth = 35.833333333333336
